# Update "想去人数" (interested-people count) figures in column F
# for the two worksheets that share this data: "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 249
    $ws.Range("F7").Value = 6822
    $ws.Range("F16").Value = 234
    $ws.Range("F17").Value = 593
}
